# Generate Report for Handback
# - Updates the "Status" text to reflect a completed handback
# - Records the Latest Target File / Latest Handback File / Latest Handback DateTime
#   for both the zh-cn and de-de localization sheets
# - Adds hyperlinks on the newly populated "Latest Target File" cells
# - Widens a few columns so the new/longer content is readable

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Latest Target File (J) / Latest Handback File (K) / Latest Handback
#    DateTime (L) for each row of the zh-cn and de-de sheets
# ---------------------------------------------------------------------------
$wsZhCn.Range("J2").Value = "de898377-1a40-4256-a36e-cf8b89338f72.md"
$wsZhCn.Range("K2").Value = "de898377-1a40-4256-a36e-cf8b89338f72.2efd77531f499766f4c2227927c500e5fb63c97f.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-02-17 09:57:02"

$wsZhCn.Range("J3").Value = "fc27428e-122e-437e-b182-a7091826a6e4.md"
$wsZhCn.Range("K3").Value = "fc27428e-122e-437e-b182-a7091826a6e4.e81b93e2e76f7862521cedf77167183d531b2600.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-02-17 09:57:02"

$wsDeDe.Range("J2").Value = "de898377-1a40-4256-a36e-cf8b89338f72.md"
$wsDeDe.Range("K2").Value = "de898377-1a40-4256-a36e-cf8b89338f72.2efd77531f499766f4c2227927c500e5fb63c97f.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-02-17 09:57:26"

$wsDeDe.Range("J3").Value = "fc27428e-122e-437e-b182-a7091826a6e4.md"
$wsDeDe.Range("K3").Value = "fc27428e-122e-437e-b182-a7091826a6e4.e81b93e2e76f7862521cedf77167183d531b2600.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-02-17 09:57:26"

# ---------------------------------------------------------------------------
# 3. Give the new "Latest Target File" cells the same hyperlink look
#    (single underline, cornflower-blue font) as the existing hyperlink
#    cells in column A, and (re)create all the hyperlinks so the new ones
#    line up with the existing "<file>.md" links.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Range("J2:J3").Font.Underline = 2
    $ws.Range("J2:J3").Font.Color = 15570276

    $deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de26946227fbf20e045aa7e5d7a6aa9c6044c8/e2e/de898377-1a40-4256-a36e-cf8b89338f72.md"
    $fcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de26946227fbf20e045aa7e5d7a6aa9c6044c8/e2e/fc27428e-122e-437e-b182-a7091826a6e4.md"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $deUrl, "", "", "de898377-1a40-4256-a36e-cf8b89338f72.md")
    $ws.Hyperlinks.Add($ws.Range("J2"), $deUrl, "", "", "de898377-1a40-4256-a36e-cf8b89338f72.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $fcUrl, "", "", "fc27428e-122e-437e-b182-a7091826a6e4.md")
    $ws.Hyperlinks.Add($ws.Range("J3"), $fcUrl, "", "", "fc27428e-122e-437e-b182-a7091826a6e4.md")
}

# ---------------------------------------------------------------------------
# 4. Widen columns to fit the new/longer content.
#    ColumnWidth is specified in characters; this engine stores width in the
#    workbook as (characters + 5/6), so we subtract 5/6 from the desired
#    stored width before assigning it.
# ---------------------------------------------------------------------------
$padding = 5 / 6

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311 - $padding
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311 - $padding

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth = 29.9777050018311 - $padding
    $ws.Columns.Item(10).ColumnWidth = 40 - $padding
    $ws.Columns.Item(11).ColumnWidth = 40 - $padding
}

Write-Host "Handback report generated"
